$wb = $excel.ActiveWorkbook

# The sheet "Foxitsoftware Reader" is the 4th sheet (sheetId=4, rId4).
$ws = $wb.Worksheets.Item(4)

# 1) Rename the sheet.
$ws.Name = "Foxitsoftware Phantompdf"

# 2) Delete three obsolete rows (bottom-to-top so row indices of the
#    still-to-be-deleted rows remain valid): original rows 13, 11, 9
#    (CVE-2018-17617, CVE-2018-17619, CVE-2018-17621). Everything below
#    shifts up, which also accounts for the dimension shrinking from
#    F99 to F96.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(9).Delete()

# 3) Replace the contents of row 5 with the new CVE record.
$ws.Range("A5").Value = "CVE-2018-17706"
$ws.Range("B5").Value = "Remote"
$ws.Range("C5").Value = "6.8"
$ws.Range("D5").Value = "2018-10-29"
$ws.Range("E5").Value = "This vulnerability allows remote attackers to execute arbitrary code on vulnerable installations of Foxit PhantomPDF Phantom PDF 9.1.5096. User interaction is required to exploit this vulnerability in that the target must visit a malicious page or open a malicious file. The specific flaw exists within fxhtml2pdf. The issue results from the lack of proper validation of user-supplied data, which can result in a memory access past the end of an allocated buffer. An attacker can leverage this vulnerability to execute code under the context of the current process. Was ZDI-CAN-6230."
$ws.Range("F5").Value = "https://www.cvedetails.com/cve/CVE-2018-17706/"
